$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns are treated as text so values
# like "7.30" or "29.40" are not coerced into numbers and lose formatting.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '70.835.86'
$ws.Range("E2").Value = '  -0.35%  '

$ws.Range("D3").Value = '3.841.36'
$ws.Range("E3").Value = '  +0.92%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = '696.11'
$ws.Range("E5").Value = '  -1.34%  '

$ws.Range("D6").Value = '171.65'
$ws.Range("E6").Value = '  -1.37%  '

$ws.Range("D7").Value = '3.840.11'
$ws.Range("E7").Value = '  +0.90%  '

$ws.Range("E9").Value = '  -0.69%  '

$ws.Range("E10").Value = '  -1.50%  '

$ws.Range("D11").Value = '7.30'
$ws.Range("E11").Value = '  -0.84%  '

$ws.Range("E12").Value = '  -1.32%  '

$ws.Range("E13").Value = '  -0.84%  '

$ws.Range("D14").Value = '36.07'
$ws.Range("E14").Value = '  -1.21%  '

$ws.Range("D15").Value = '4.492.64'
$ws.Range("E15").Value = '  +1.02%  '

$ws.Range("D16").Value = '3.849.66'
$ws.Range("E16").Value = '  +1.22%  '

$ws.Range("D17").Value = '70.866.97'
$ws.Range("E17").Value = '  -0.31%  '

$ws.Range("E18").Value = '  -1.48%  '

$ws.Range("E19").Value = '  +0.43%  '

$ws.Range("D20").Value = '17.32'
$ws.Range("E20").Value = '  -3.99%  '

$ws.Range("E21").Value = '  -4.83%  '

$ws.Range("D22").Value = '494.18'
$ws.Range("E22").Value = '  +2.07%  '

$ws.Range("D23").Value = '0.716'
$ws.Range("E23").Value = '  -0.50%  '

$ws.Range("D24").Value = '84.56'
$ws.Range("E24").Value = '  +0.61%  '

$ws.Range("E25").Value = '  +0.76%  '

$ws.Range("B26").Value = 'RenderToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D26").Value = '10.55'
$ws.Range("E26").Value = '  -1.15%  '

$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").Value = '12.14'
$ws.Range("E27").Value = '  -3.67%  '

$ws.Range("E28").Value = '  -4.05%  '

$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.06%  '

$ws.Range("E30").Value = '  -0.63%  '

$ws.Range("E31").Value = '  -2.44%  '

$ws.Range("D32").Value = '2.25'
$ws.Range("E32").Value = '  -2.82%  '

$ws.Range("D33").Value = '29.40'
$ws.Range("E33").Value = '  -0.98%  '

$ws.Range("E34").Value = '  -0.58%  '

$ws.Range("D35").Value = '3.798.71'
$ws.Range("E35").Value = '  +1.13%  '

$ws.Range("D36").Value = '9.13'
$ws.Range("E36").Value = '  -2.13%  '

$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  +0.05%  '

$ws.Range("E38").Value = '  -0.87%  '

$ws.Range("D39").Value = '2.38'
$ws.Range("E39").Value = '  +5.66%  '

$ws.Range("E40").Value = '  +6.83%  '

$ws.Range("E41").Value = '  -1.07%  '

$ws.Range("E42").Value = '  -5.53%  '

$ws.Range("E43").Value = '  +0.03%  '

$ws.Range("E44").Value = '  +0.15%  '

$ws.Range("B45").Value = 'FLOKI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D45").Value = '0.000311'
$ws.Range("E45").Value = '  -6.16%  '

$ws.Range("B46").Value = 'Monero'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D46").Value = '163.36'
$ws.Range("E46").Value = '  +1.45%  '

$ws.Range("D47").Value = '48.74'
$ws.Range("E47").Value = '  -1.47%  '

$ws.Range("D48").Value = '0.298'
$ws.Range("E48").Value = '  -1.24%  '

$ws.Range("D49").Value = '8.61'
$ws.Range("E49").Value = '  +0.28%  '

$ws.Range("D50").Value = '43.08'
$ws.Range("E50").Value = '  -5.82%  '

$ws.Range("D51").Value = '404.16'
$ws.Range("E51").Value = '  -0.25%  '
